$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Activate() | Out-Null

# Insert a new blank column before column N (the existing "Late"/"Outstanding"
# columns shift right by one, from N/P to O/Q) - this mirrors a manual
# "Insert Column" in Excel, which carries the left neighbour's width onto
# the freshly inserted column.
$leftWidth = $ws.Columns("M").ColumnWidth
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $leftWidth

# Leave the selection on the newly active sheet parked at S7, matching the
# state captured after the edit.
$ws.Range("S7").Select() | Out-Null
